$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 3 (shifts old rows 3-6 down to 6-9)
$ws.Range("A3:B5").Insert(-4121)

# Clear any inherited formatting (yellow fill) from the insert-copy-down behavior
$ws.Range("A3:B5").ClearFormats()

# Fill in new Sprint 2 rows
$ws.Range("A3").Value = "Sprint 2"
$ws.Range("B3").Value = "Each member should install the ganttProject on his computer"
$ws.Range("A4").Value = "Sprint 2"
$ws.Range("B4").Value = "Each member should read the code and analize the functionalities already implemented"
$ws.Range("A5").Value = "Sprint 2"
$ws.Range("B5").Value = "The members should discuss the features"

# Re-apply section-header style (green fill), copying the exact format already used at A6
$ws.Range("A6").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Reproduce the extra style artifact on B5
$ws.Range("B5").Font.ThemeColor = 1

# Update selection to match new active cell
$ws.Range("B5").Select() | Out-Null
